$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.014855
$ws.Range("H2").Value = 0.044565
$ws.Range("I2").Value = 0.5536850213696453
$ws.Range("J2").Value = 0.5536850213696451
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005328
$ws.Range("N2").Value = 0.015984
$ws.Range("O2").Value = 0.001689940172269439
$ws.Range("P2").Value = 0.001689940172269439
$ws.Range("Q2").Value = [double]"7.914744000000001E-05"
$ws.Range("R2").Value = 0.0007123269600000001
$ws.Range("S2").Value = 0.0009356945603964266
$ws.Range("T2").Value = 0.0009356945603964264
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.014855
$ws.Range("H3").Value = 0.044565
$ws.Range("I3").Value = 0.5536850213696453
$ws.Range("J3").Value = 0.5536850213696451
$ws.Range("O3").Value = 0.9983100598277306
$ws.Range("P3").Value = 0.9983100598277306
$ws.Range("Q3").Value = 0.04675531528166667
$ws.Range("R3").Value = 0.4207978375350001
$ws.Range("S3").Value = 0.5527493268092488
$ws.Range("T3").Value = 0.5527493268092487
$ws.Range("G4").Value = 0.01197433333333333
$ws.Range("I4").Value = 0.4463149786303549
$ws.Range("J4").Value = 0.4463149786303548
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.005328
$ws.Range("N4").Value = 0.015984
$ws.Range("O4").Value = 0.001689940172269439
$ws.Range("P4").Value = 0.001689940172269439
$ws.Range("Q4").Value = [double]"6.379924799999999E-05"
$ws.Range("R4").Value = 0.000574193232
$ws.Range("S4").Value = 0.000754245611873013
$ws.Range("T4").Value = 0.0007542456118730129
$ws.Range("G5").Value = 0.01197433333333333
$ws.Range("I5").Value = 0.4463149786303549
$ws.Range("J5").Value = 0.4463149786303548
$ws.Range("O5").Value = 0.9983100598277306
$ws.Range("P5").Value = 0.9983100598277306
$ws.Range("R5").Value = 0.339197143897
$ws.Range("S5").Value = 0.4455607330184819
$ws.Range("T5").Value = 0.4455607330184818
